# Generate Report for Handback
# Insert a new row for file "910bdb67-0a00-40aa-bd98-949a431a4405.md"
# (handed back: in sync with en-US) between the existing "84f8000b..."
# row and the existing "de0b7adf..." row, on all three sheets
# (Overview, zh-cn, de-de), and keep the tables / dimensions / hyperlinks
# in sync.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows(3).Insert()

$ws.Range("A3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.md"
$ws.Range("B3").Value = "'e2e\910bdb67-0a00-40aa-bd98-949a431a4405.md"
$ws.Range("C3").Value = "'.md"
$ws.Range("E3").Value = "'Handed back: in sync with en-US"
$ws.Range("F3").Value = "'Handed back: in sync with en-US"
$ws.Range("G3").Value = "'2016-08-26 16:45:26"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G4"))

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b7c200ac9842fbee428563245f85417fcdf853e/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md", "", "", "e2e\84f8000b-cf6b-4b2a-8656-2428c2200111.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08/e2e/910bdb67-0a00-40aa-bd98-949a431a4405.md", "", "", "e2e\910bdb67-0a00-40aa-bd98-949a431a4405.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b31178e617427d1b15fa2f67ded880d9448843ae/e2e/de0b7adf-3cc1-489c-89d1-d3a621299d05.md", "", "", "e2e\de0b7adf-3cc1-489c-89d1-d3a621299d05.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows(3).Insert()

$ws.Range("A3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.md"
$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'Handed back: in sync with en-US"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08.zh-cn.xlf"
$ws.Range("H3").Value = "'2016-08-26 16:45:22"
$ws.Range("I3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.md"
$ws.Range("J3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08.zh-cn.xlf"
$ws.Range("K3").Value = "'2016-08-26 16:45:38"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P4"))

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b7c200ac9842fbee428563245f85417fcdf853e/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md", "", "", "84f8000b-cf6b-4b2a-8656-2428c2200111.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/54004e9e964e93fbf39575ecad02cae92a29df44/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md", "", "", "84f8000b-cf6b-4b2a-8656-2428c2200111.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08/e2e/910bdb67-0a00-40aa-bd98-949a431a4405.md", "", "", "910bdb67-0a00-40aa-bd98-949a431a4405.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08/e2e/910bdb67-0a00-40aa-bd98-949a431a4405.md", "", "", "910bdb67-0a00-40aa-bd98-949a431a4405.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b31178e617427d1b15fa2f67ded880d9448843ae/e2e/de0b7adf-3cc1-489c-89d1-d3a621299d05.md", "", "", "de0b7adf-3cc1-489c-89d1-d3a621299d05.md")
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4e60d5ef2d63ca777440ed615f5691f67caf2487/e2e/de0b7adf-3cc1-489c-89d1-d3a621299d05.md", "", "", "de0b7adf-3cc1-489c-89d1-d3a621299d05.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows(3).Insert()

$ws.Range("A3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.md"
$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'Handed back: in sync with en-US"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08.de-de.xlf"
$ws.Range("H3").Value = "'2016-08-26 16:45:26"
$ws.Range("I3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.md"
$ws.Range("J3").Value = "'910bdb67-0a00-40aa-bd98-949a431a4405.9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08.de-de.xlf"
$ws.Range("K3").Value = "'2016-08-26 16:45:44"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P4"))

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b7c200ac9842fbee428563245f85417fcdf853e/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md", "", "", "84f8000b-cf6b-4b2a-8656-2428c2200111.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6b996c91f7550691d093844f49bbe4eb92b9f5f1/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md", "", "", "84f8000b-cf6b-4b2a-8656-2428c2200111.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08/e2e/910bdb67-0a00-40aa-bd98-949a431a4405.md", "", "", "910bdb67-0a00-40aa-bd98-949a431a4405.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9e53fc9d3b5e171d8041b1ac0712f37ef70f2c08/e2e/910bdb67-0a00-40aa-bd98-949a431a4405.md", "", "", "910bdb67-0a00-40aa-bd98-949a431a4405.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b31178e617427d1b15fa2f67ded880d9448843ae/e2e/de0b7adf-3cc1-489c-89d1-d3a621299d05.md", "", "", "de0b7adf-3cc1-489c-89d1-d3a621299d05.md")
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/41d27a1b7f2c144c9aa184e7e8bf9ca2e08dbddc/e2e/de0b7adf-3cc1-489c-89d1-d3a621299d05.md", "", "", "de0b7adf-3cc1-489c-89d1-d3a621299d05.md")
